# Update data rows 2-5 with new dataset values (PLoP), and widen selected columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: overwrite row data for rows 2-5 (A:AH) ---
$row2 = @(45039.50694444445,24.502,17.071,4.266,51.765,42.826,19.282,64.646,29.668,12.708,19.604,20.148,21.386,6.157,19.174,27.062,15.864,3.991,2.607,284.514,53.424,17.698,35.645,18.491,2.496,31.956,15.633,13.958,16.32,21.117,3.641,57.291,9.937,22.127)
$arr2 = New-Object 'object[,]' 1,34
for ($i = 0; $i -lt 34; $i++) { $arr2[0,$i] = $row2[$i] }
$ws.Range("A2:AH2").Value = $arr2

$row3 = @(45039.51388888889,1.441,0.416,1.313,2.538,2.157,1.14,12.856,1.745,0.658,0.962,0.823,0.896,0.402,1.128,1.646,1.181,1.535,0.619,10.02,3.805,1.041,2.421,1.272,0.131,5.24,0.92,1.069,1.181,0.868,1.294,12.108,0.467,1.316)
$arr3 = New-Object 'object[,]' 1,34
for ($i = 0; $i -lt 34; $i++) { $arr3[0,$i] = $row3[$i] }
$ws.Range("A3:AH3").Value = $arr3

$row4 = @(45039.52083333334,9.128,6.414,1.079,19.467,16.155,7.183,25.811,11.053,4.722,7.258,7.705,8.168,2.29,7.143,10.004,6.154,1.038,0.594,101.354,19.805,6.594,13.082,7.122,0.92,12.157,5.824,5.293,6.192,8.138,0.784,22.804,3.654,8.242)
$arr4 = New-Object 'object[,]' 1,34
for ($i = 0; $i -lt 34; $i++) { $arr4[0,$i] = $row4[$i] }
$ws.Range("A4:AH4").Value = $arr4

$row5 = @(45039.52777777778,5.77,4,0.75,12.25,10.18,4.54,19.11,6.98,2.99,4.56,4.84,5.14,1.45,4.51,6.35,3.95,0.75,0.39,61.35,12.7,4.16,8.38,4.56,0.57,8.91,3.68,3.37,3.94,5.14,0.56,17.29,2.28,5.21)
$arr5 = New-Object 'object[,]' 1,34
for ($i = 0; $i -lt 34; $i++) { $arr5[0,$i] = $row5[$i] }
$ws.Range("A5:AH5").Value = $arr5

# --- Step 2: remove old row 6 (dataset now has only 4 data rows) ---
$ws.Rows.Item(6).Delete()

# --- Step 3: widen columns B,C,G,J,K,L,M,O,Q,V,X,AA,AB,AD,AH from 7 to 8 ---
# ColumnWidth (chars) -> stored OOXML width has a fixed +0.8333333333333333 offset
# in this engine, so to land on stored width 8 we set ColumnWidth = 8 - 0.8333333333333333
$colsToWiden = @(2,3,7,10,11,12,13,15,17,22,24,27,28,30,34)
foreach ($c in $colsToWiden) {
    $ws.Columns.Item($c).ColumnWidth = 7.166666666666667
}
